$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / shared text updates ---
$ws.Range("M6").Value = "Edward A. Caban"
$ws.Range("A8").Value = "Volume 30   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# --- Crime data table updates (rows 14-30) ---
# NOTE: some cells flip between a numeric value and the sheet's
# "N/A" text placeholders (shared strings "0" / "***.*"). For those cells we
# first force the desired value/type, then copy cell formatting from a stable
# donor cell of the same target style so the style index matches the sheet's
# convention (style 14 = text placeholder look, style 15 = numeric look).

# Row 14
$ws.Range("M14").Value = -33.333333333333

# Row 15
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "***.*"
$ws.Range("L15").Value = 23.529411764705
$ws.Range("M15").Value = 110

# Row 16
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -4.347826086956
$ws.Range("I16").Value = 142
$ws.Range("J16").Value = 121
$ws.Range("K16").Value = 17.355371900826
$ws.Range("L16").Value = 30.275229357798
$ws.Range("M16").Value = -10.126582278481
$ws.Range("N16").Value = -84.812834224598

# Row 17
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = 11.111111111111
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 43
$ws.Range("H17").Value = -25.581395348837
$ws.Range("I17").Value = 246
$ws.Range("J17").Value = 274
$ws.Range("K17").Value = -10.218978102189
$ws.Range("L17").Value = 12.844036697247
$ws.Range("M17").Value = 66.216216216216
$ws.Range("N17").Value = -39.853300733496

# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -55.555555555555
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 37
$ws.Range("H18").Value = -70.27027027027
$ws.Range("I18").Value = 112
$ws.Range("J18").Value = 153
$ws.Range("K18").Value = -26.797385620915
$ws.Range("L18").Value = -20.567375886524
$ws.Range("M18").Value = -54.655870445344
$ws.Range("N18").Value = -89.532710280373

# Row 19
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 54.545454545454
$ws.Range("F19").Value = 73
$ws.Range("G19").Value = 67
$ws.Range("H19").Value = 8.955223880597
$ws.Range("I19").Value = 472
$ws.Range("J19").Value = 432
$ws.Range("K19").Value = 9.259259259259
$ws.Range("L19").Value = 38.41642228739
$ws.Range("M19").Value = 37.209302325581
$ws.Range("N19").Value = -17.91304347826

# Row 20
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 36.363636363636
$ws.Range("F20").Value = 46
$ws.Range("G20").Value = 43
$ws.Range("H20").Value = 6.976744186046
$ws.Range("I20").Value = 198
$ws.Range("J20").Value = 150
$ws.Range("K20").Value = 32
$ws.Range("L20").Value = 90.384615384615
$ws.Range("M20").Value = 32.885906040268
$ws.Range("N20").Value = -85.662563359884

# Row 21
$ws.Range("C21").Value = 54
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 22.727272727272
$ws.Range("F21").Value = 186
$ws.Range("G21").Value = 213
$ws.Range("H21").Value = -12.676056338028
$ws.Range("I21").Value = 1195
$ws.Range("J21").Value = 1149
$ws.Range("K21").Value = 4.003481288076
$ws.Range("L21").Value = 27.534685165421
$ws.Range("M21").Value = 12.523540489642
$ws.Range("N21").Value = -72.871736662883

# Row 22
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("J22").Value = 12
$ws.Range("K22").Value = 66.666666666666

# Row 23
$ws.Range("C23").Value = 6
$ws.Range("E23").Value = 20
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = -17.647058823529
$ws.Range("I23").Value = 122
$ws.Range("J23").Value = 114
$ws.Range("K23").Value = 7.017543859649
$ws.Range("L23").Value = 23.232323232323
$ws.Range("M23").Value = 35.555555555555

# Row 24
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 67
$ws.Range("E24").Value = -46.268656716417
$ws.Range("F24").Value = 153
$ws.Range("G24").Value = 241
$ws.Range("H24").Value = -36.514522821576
$ws.Range("I24").Value = 1205
$ws.Range("J24").Value = 1363
$ws.Range("K24").Value = -11.592076302274
$ws.Range("L24").Value = -0.659521846661
$ws.Range("M24").Value = 56.493506493506

# Row 25
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -29.411764705882
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 90
$ws.Range("H25").Value = -33.333333333333
$ws.Range("I25").Value = 439
$ws.Range("J25").Value = 517
$ws.Range("K25").Value = -15.087040618955
$ws.Range("L25").Value = 5.023923444976
$ws.Range("M25").Value = -4.357298474945

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 29
$ws.Range("J26").Value = 25
$ws.Range("K26").Value = 16
$ws.Range("L26").Value = 3.571428571428

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 12
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 53
$ws.Range("K27").Value = 35.897435897435
$ws.Range("L27").Value = -18.461538461538

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 17
$ws.Range("K28").Value = 183.333333333333
$ws.Range("L28").Value = 13.333333333333
$ws.Range("M28").Value = 30.76923076923
$ws.Range("N28").Value = 0

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 16
$ws.Range("K29").Value = 166.666666666667
$ws.Range("L29").Value = 33.333333333333
$ws.Range("M29").Value = 45.454545454545
$ws.Range("N29").Value = 0

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"

# --- Fix up style indices for the placeholder/value type-flip cells so they
# match the sheet's existing text-placeholder (style 14) / numeric (style 15)
# look, by copying formatting from stable donor cells. ---
$ws.Range("C14").Copy() | Out-Null
$ws.Range("G15").PasteSpecial(-4122) | Out-Null
$ws.Range("H15").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("F14").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
